$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.831.71'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '1.905.14'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.06'
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5037'
$ws.Range("E7").Value = '  +4.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3810'
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9065'
$ws.Range("E10").Value = '  -2.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.89'
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").Value = '1.967.65'
$ws.Range("E12").Value = '  +2.92%  '
$ws.Range("E13").Value = '  -1.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.495'
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.80'
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008709'
$ws.Range("E17").Value = '  -1.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").Value = '27.864.49'
$ws.Range("E19").Value = '  -0.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.61'
$ws.Range("E20").Value = '  -1.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.160'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.578'
$ws.Range("E23").Value = '  -0.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '153.70'
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.878'
$ws.Range("E25").Value = '  -2.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.223'
$ws.Range("E26").Value = '  +4.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.37'
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.34'
$ws.Range("E28").Value = '  -1.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.906'
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09028'
$ws.Range("E30").Value = '  +1.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.213'
$ws.Range("E31").Value = '  -2.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.220'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.707'
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7675'
$ws.Range("E34").Value = '  -0.79%  '
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.507'
$ws.Range("E36").Value = '  -5.07%  '
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5531'
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.014'
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05262'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.886'
$ws.Range("E41").Value = '  -2.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1510'
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '111.33'
$ws.Range("E44").Value = '  +3.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.60'
$ws.Range("E45").Value = '  -0.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4807'
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.629'
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.51'
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06065'
$ws.Range("E50").Value = '  -0.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9031'
$ws.Range("E51").Value = '  +0.43%  '
